$wb = $excel.ActiveWorkbook

# --- Sheet "Top Gainers": rows 36-76 shift up by one (new row added at bottom, oldest row dropped) ---
$wsGainers = $wb.Worksheets.Item("Top Gainers")

function Set-GainerRow($r, $b, $c, $d, $e) {
    $wsGainers.Cells.Item($r, 2).Value = $b
    $wsGainers.Cells.Item($r, 3).Value = $c
    $wsGainers.Cells.Item($r, 4).Value = $d
    $wsGainers.Cells.Item($r, 5).Value = $e
}

Set-GainerRow 36 "PROZONER" 4.9921 15.7468 36.095
Set-GainerRow 37 "STALLION" 4.9914 -5.2229 21.4391
Set-GainerRow 38 "SURYAROSNI" 4.9831 11.386 3.0213
Set-GainerRow 39 "DATAMATICS" 4.9005 7.3152 15.7298
Set-GainerRow 40 "UTKARSHBNK" 4.8768 -5.8959 -2.6215
Set-GainerRow 41 "FILATEX" 4.8689 10.274 26.0027
Set-GainerRow 42 "INDOTHAI" 4.8064 4.5349 43.748
Set-GainerRow 43 "SANDUMA" 4.593 2.1405 30.2813
Set-GainerRow 44 "LLOYDSENT" 4.5646 1.8339 11.234
Set-GainerRow 45 "STAR" 4.5025 4.4319 3.662
Set-GainerRow 46 "RECLTD" 4.4992 3.4756 3.4062
Set-GainerRow 47 "NBCC" 4.4511 3.1605 7.6018
Set-GainerRow 48 "GPPL" 4.4154 3.4073 5.0497
Set-GainerRow 49 "BIL" 4.3654 9.122199999999999 -0.3203
Set-GainerRow 50 "HUDCO" 4.3201 3.8924 5.3884
Set-GainerRow 51 "SGMART" 4.2736 8.258900000000001 2.5381
Set-GainerRow 52 "MRPL" 4.2642 9.7103 20.0542
Set-GainerRow 53 "JKIL" 4.1372 2.9463 1.7584
Set-GainerRow 54 "SAMBHV" 4.1349 2.624 5.167
Set-GainerRow 55 "SAPPHIRE" 4.1265 1.7633 -0.7999000000000001
Set-GainerRow 56 "PVRINOX" 4.1118 6.2102 14.707
Set-GainerRow 57 "KERNEX" 3.9981 7.4592 27.1054
Set-GainerRow 58 "SUNFLAG" 3.997 4.333 4.6312
Set-GainerRow 59 "CMSINFO" 3.9096 2.6872 2.8935
Set-GainerRow 60 "GMBREW" 3.8999 -0.53 79.029
Set-GainerRow 61 "APARINDS" 3.8924 8.3414 15.5876
Set-GainerRow 62 "HITECHGEAR" 3.8587 1.1486 9.9254
Set-GainerRow 63 "NPST" 3.8509 -2.0059 -3.5057
Set-GainerRow 64 "ORIENTTECH" 3.827 0.5247000000000001 32.6784
Set-GainerRow 65 "ICRA" 3.7985 4.4793 2.8828
Set-GainerRow 66 "SALASAR" 3.7935 4.7872 11.0485
Set-GainerRow 67 "DCW" 3.7544 2.3219 -3.9753
Set-GainerRow 68 "RHETAN" 3.754 4.178 6.549
Set-GainerRow 69 "HINDPETRO" 3.6935 6.9335 5.7397
Set-GainerRow 70 "BHARTIHEXA" 3.6718 7.0877 15.3332
Set-GainerRow 71 "HLEGLAS" 3.659 8.115500000000001 27.1239
Set-GainerRow 72 "RHIM" 3.6544 3.2276 5.1826
Set-GainerRow 73 "SHK" 3.6347 2.388 -1.932
Set-GainerRow 74 "BCLIND" 3.6271 2.2945 0.1728
Set-GainerRow 75 "MUKANDLTD" 3.6133 11.9685 9.550800000000001
Set-GainerRow 76 "CGPOWER" 3.6125 3.4192 1.0325

# --- Sheet "Top Losers": update individual Weekly values ---
$wsLosers = $wb.Worksheets.Item("Top Losers")
$wsLosers.Cells.Item(18, 4).Value = -0.062
$wsLosers.Cells.Item(48, 4).Value = 0.05
$wsLosers.Cells.Item(56, 4).Value = 3.7771
